# Update workbook/sheet name and data to reflect new "through" date (2022-06-08 -> 2022-06-09)
# and updated values for the month of July (row 7) and the Total row (row 14), column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab
$ws.Name = "Through 2022-06-09"

# Update the header label in column I (year-to-date column header)
$ws.Range("I1").Value = "2022 (through 06-09)"

# Update the data values
$ws.Range("I7").Value = 34
$ws.Range("I14").Value = 697
